$wb = $excel.ActiveWorkbook

# --- Full Message Examples sheet: update the "CAN Command Example" and
# related example rows to reflect commands now taking 4 U8 parameters
# instead of 1. ---
$ws = $wb.Worksheets.Item("Full Message Examples")

# Request Example (row 9-11): parameter id 0x2 -> 0x1
$ws.Range("E10").Value = "0x1"

# CAN Command Example (row 13-15)
$ws.Range("E14").Value = "0x0"    # command id 0x1 -> 0x0
$ws.Range("H14").Value = "0x5"    # DLC 0x2 -> 0x5 (1 cmd byte + 4 param bytes)
$ws.Range("L14").Value = "0x00"   # new param byte 2
$ws.Range("M14").Value = "0x00"   # new param byte 3
$ws.Range("N14").Value = "0x00"   # new param byte 4

$ws.Range("L15").Value = "Parameter 1"
$ws.Range("M15").Value = "Parameter 2"
$ws.Range("N15").Value = "Parameter 3"

# CAN Error Example (row 17-19): parameter id 0x2 -> 0x1
$ws.Range("E18").Value = "0x1"

# --- View / selection state restored to match author's workbook ---
$ws1 = $wb.Worksheets.Item("ID Example")
$ws1.Activate()
$ws1.Range("O10").Select()

$ws.Activate()
$ws.Range("H16").Select()
